$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concise")
$ws.Range("C3").Value = 233
